$wb = $excel.ActiveWorkbook

# --- ALC sheet: scheduled market-price refresh ---
$ws = $wb.Worksheets.Item("ALC")

# Row 40
$ws.Range("H40").Value = 1484.2693
$ws.Range("I40").Value = 1480.9524
$ws.Range("J40").Value = 1498.2
$ws.Range("K40").Value = 1480.9524
$ws.Range("L40").Value = 1498.2
$ws.Range("M40").Value = -1305.9524
$ws.Range("N40").Value = -1848.2

# Row 76
$ws.Range("H76").Value = 3471.2942
$ws.Range("I76").Value = 3125.375
$ws.Range("K76").Value = 3125.375
$ws.Range("M76").Value = -2810.375

# Row 79
$ws.Range("H79").Value = 3471.2942
$ws.Range("I79").Value = 3125.375
$ws.Range("K79").Value = 3125.375
$ws.Range("M79").Value = -2033.375

# Row 86
$ws.Range("H86").Value = 8138.75
$ws.Range("I86").Value = 2380.4
$ws.Range("J86").Value = 17736
$ws.Range("K86").Value = 2380.4
$ws.Range("L86").Value = 17736
$ws.Range("M86").Value = -1257.4
$ws.Range("N86").Value = -19982

# Row 89
$ws.Range("H89").Value = 8138.75
$ws.Range("I89").Value = 2380.4
$ws.Range("J89").Value = 17736
$ws.Range("K89").Value = 11902
$ws.Range("L89").Value = 88680
$ws.Range("M89").Value = -6286
$ws.Range("N89").Value = -99912

# Row 98
$ws.Range("H98").Value = 990.9048
$ws.Range("I98").Value = 874.1579
$ws.Range("K98").Value = 874.1579
$ws.Range("M98").Value = 623.8421

# Row 122
$ws.Range("H122").Value = 990.9048
$ws.Range("I122").Value = 874.1579
$ws.Range("K122").Value = 2622.4737
$ws.Range("M122").Value = -172.4737

# Row 132
$ws.Range("H132").Value = 15849.714
$ws.Range("I132").Value = 18227
$ws.Range("J132").Value = 1586
$ws.Range("K132").Value = 54681
$ws.Range("L132").Value = 4758
$ws.Range("M132").Value = -52151
$ws.Range("N132").Value = -9818

# Row 138
$ws.Range("H138").Value = 2806.6191
$ws.Range("I138").Value = 3380.6365
$ws.Range("K138").Value = 10141.9095
$ws.Range("M138").Value = -5001.9095

# Row 141
$ws.Range("H141").Value = 5533.579
$ws.Range("I141").Value = 5040.1875
$ws.Range("J141").Value = 8165
$ws.Range("K141").Value = 15120.5625
$ws.Range("L141").Value = 24495
$ws.Range("M141").Value = -9940.5625
$ws.Range("N141").Value = -34855

# --- ARM sheet: scheduled market-price refresh ---
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 251560.56
$ws.Range("I32").Value = 298935.88
$ws.Range("K32").Value = 298935.88
$ws.Range("M32").Value = -298648.88

# Row 45
$ws.Range("H45").Value = 3307.3635
$ws.Range("I45").Value = 3425
$ws.Range("K45").Value = 3425
$ws.Range("M45").Value = -3048

# Row 132
$ws.Range("H132").Value = 5475.754
$ws.Range("I132").Value = 3920.1667
$ws.Range("K132").Value = 11760.5001
$ws.Range("M132").Value = -9230.500100000001

# --- BSM sheet: scheduled market-price refresh ---
$ws = $wb.Worksheets.Item("BSM")

# Row 80
$ws.Range("H80").Value = 650.64
$ws.Range("I80").Value = 682.7273
$ws.Range("J80").Value = 625.4286
$ws.Range("K80").Value = 682.7273
$ws.Range("L80").Value = 625.4286
$ws.Range("M80").Value = 315.2727
$ws.Range("N80").Value = -2621.4286

# Row 83
$ws.Range("H83").Value = 650.64
$ws.Range("I83").Value = 682.7273
$ws.Range("J83").Value = 625.4286
$ws.Range("K83").Value = 3413.6365
$ws.Range("L83").Value = 3127.143
$ws.Range("M83").Value = 1578.3635
$ws.Range("N83").Value = -13111.143

# Row 86
$ws.Range("H86").Value = 1961.8
$ws.Range("I86").Value = 1899.2
$ws.Range("J86").Value = 2149.6
$ws.Range("K86").Value = 1899.2
$ws.Range("L86").Value = 2149.6
$ws.Range("M86").Value = -776.2
$ws.Range("N86").Value = -4395.6

# Row 89
$ws.Range("H89").Value = 1961.8
$ws.Range("I89").Value = 1899.2
$ws.Range("J89").Value = 2149.6
$ws.Range("K89").Value = 9496
$ws.Range("L89").Value = 10748
$ws.Range("M89").Value = -3880
$ws.Range("N89").Value = -21980

# Row 105
$ws.Range("H105").Value = 1835
$ws.Range("I105").Value = 1646.5358
$ws.Range("K105").Value = 1646.5358
$ws.Range("M105").Value = 100.4641999999999

# --- CRP sheet: scheduled market-price refresh ---
$ws = $wb.Worksheets.Item("CRP")

# Row 70
$ws.Range("H70").Value = 37600
$ws.Range("J70").Value = 37600
$ws.Range("L70").Value = 37600
$ws.Range("N70").Value = -38230

# Row 73
$ws.Range("H73").Value = 37600
$ws.Range("J73").Value = 37600
$ws.Range("L73").Value = 37600
$ws.Range("N73").Value = -39784

# Row 99
$ws.Range("H99").Value = 12010.857
$ws.Range("I99").Value = 24317.889
$ws.Range("K99").Value = 24317.889
$ws.Range("M99").Value = -22819.889

# Row 105
$ws.Range("H105").Value = 13469.625
$ws.Range("I105").Value = 13469.625
$ws.Range("K105").Value = 13469.625
$ws.Range("M105").Value = -11722.625

# Row 126
$ws.Range("H126").Value = 12010.857
$ws.Range("I126").Value = 24317.889
$ws.Range("K126").Value = 72953.667
$ws.Range("M126").Value = -70483.667

# --- CUL sheet: scheduled market-price refresh ---
$ws = $wb.Worksheets.Item("CUL")

# Row 137
$ws.Range("H137").Value = 6900.8335
$ws.Range("I137").Value = 1837.8182
$ws.Range("K137").Value = 5513.4546
$ws.Range("M137").Value = -413.4546

# Row 140
$ws.Range("H140").Value = 1895.9
$ws.Range("I140").Value = 1773.2222
$ws.Range("K140").Value = 5319.6666
$ws.Range("M140").Value = -139.6665999999996

# --- GSM sheet: scheduled market-price refresh ---
$ws = $wb.Worksheets.Item("GSM")

# Row 15
$ws.Range("H15").Value = 22007.334
$ws.Range("J15").Value = 19758.25
$ws.Range("L15").Value = 19758.25
$ws.Range("N15").Value = -20334.25

# Row 42
$ws.Range("H42").Value = 60087
$ws.Range("J42").Value = 59999
$ws.Range("L42").Value = 59999
$ws.Range("N42").Value = -60969

# Row 43
$ws.Range("H43").Value = 1874.6666
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0

# Row 45
$ws.Range("H45").Value = 47000
$ws.Range("J45").Value = 47000
$ws.Range("L45").Value = 47000
$ws.Range("N45").Value = -48118

# Row 81
$ws.Range("H81").Value = 22007.334
$ws.Range("J81").Value = 19758.25
$ws.Range("L81").Value = 19758.25
$ws.Range("N81").Value = -21754.25

# Row 84
$ws.Range("H84").Value = 22007.334
$ws.Range("J84").Value = 19758.25
$ws.Range("L84").Value = 59274.75
$ws.Range("N84").Value = -69258.75

# Row 86
$ws.Range("H86").Value = 9990
$ws.Range("I86").Value = 9990
$ws.Range("K86").Value = 9990
$ws.Range("M86").Value = -8804

# Row 89
$ws.Range("H89").Value = 9990
$ws.Range("I89").Value = 9990
$ws.Range("K89").Value = 29970
$ws.Range("M89").Value = -24042

# Row 115
$ws.Range("H115").Value = 60087
$ws.Range("J115").Value = 59999
$ws.Range("L115").Value = 59999
$ws.Range("N115").Value = -62349

# Row 43: N43 cleared (no longer applicable)
$ws.Range("N43").ClearContents()

# --- LTW sheet: scheduled market-price refresh ---
$ws = $wb.Worksheets.Item("LTW")

# Row 68
$ws.Range("H68").Value = 6571.2173
$ws.Range("I68").Value = 7936.6665
$ws.Range("J68").Value = 4011
$ws.Range("K68").Value = 7936.6665
$ws.Range("L68").Value = 4011
$ws.Range("M68").Value = -7187.6665
$ws.Range("N68").Value = -5509

# Row 71
$ws.Range("H71").Value = 6571.2173
$ws.Range("I71").Value = 7936.6665
$ws.Range("J71").Value = 4011
$ws.Range("K71").Value = 39683.3325
$ws.Range("L71").Value = 20055
$ws.Range("M71").Value = -35939.3325
$ws.Range("N71").Value = -27543

# Row 80
$ws.Range("H80").Value = 44500
$ws.Range("J80").Value = 44500
$ws.Range("L80").Value = 44500
$ws.Range("N80").Value = -46746

# Row 83
$ws.Range("H83").Value = 44500
$ws.Range("J83").Value = 44500
$ws.Range("L83").Value = 133500
$ws.Range("N83").Value = -144732

# Row 133
$ws.Range("H133").Value = 88777
$ws.Range("J133").Value = 88777
$ws.Range("L133").Value = 88777
$ws.Range("N133").Value = -93837

# --- WVR sheet: scheduled market-price refresh ---
$ws = $wb.Worksheets.Item("WVR")

# Row 81
$ws.Range("H81").Value = 52782.1
$ws.Range("I81").Value = 2099.0833
$ws.Range("J81").Value = 128806.625
$ws.Range("K81").Value = 4198.1666
$ws.Range("L81").Value = 257613.25
$ws.Range("M81").Value = -3137.1666
$ws.Range("N81").Value = -259735.25

# Row 84
$ws.Range("H84").Value = 52782.1
$ws.Range("I84").Value = 2099.0833
$ws.Range("J84").Value = 128806.625
$ws.Range("K84").Value = 20990.833
$ws.Range("L84").Value = 1288066.25
$ws.Range("M84").Value = -15686.833
$ws.Range("N84").Value = -1298674.25

# Row 136
$ws.Range("H136").Value = 1409.16
$ws.Range("I136").Value = 1079.5555
$ws.Range("J136").Value = 2256.7144
$ws.Range("K136").Value = 3238.6665
$ws.Range("L136").Value = 6770.1432
$ws.Range("M136").Value = -688.6664999999998
$ws.Range("N136").Value = -11870.1432
